# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.161.92"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.359.44"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'0.677"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").Value = "'238.54"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").Value = "'72.99"
$ws.Range("E7").Value = "  +11.82%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "  +19.45%  "
$ws.Range("E10").Value = "  +6.63%  "
$ws.Range("D11").Value = "'29.46"
$ws.Range("E11").Value = "  +10.78%  "
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "2.713.81"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'16.82"
$ws.Range("E15").Value = "  +7.72%  "
$ws.Range("D16").Value = "'0.905"
$ws.Range("E16").Value = "  +8.01%  "
$ws.Range("D17").Value = "2.362.52"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "44.109.66"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").Value = "'77.96"
$ws.Range("E20").Value = "  +6.05%  "
$ws.Range("E21").Value = "  +4.77%  "
$ws.Range("D22").Value = "'255.65"
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'3.78"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").Value = "'10.51"
$ws.Range("E26").Value = "  +6.79%  "
$ws.Range("D27").Value = "'2.23"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "'172.92"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").Value = "'1.59"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("E32").Value = "  +5.34%  "
$ws.Range("E33").Value = "  +4.96%  "
$ws.Range("E34").Value = "  +6.95%  "
$ws.Range("D35").Value = "'5.23"
$ws.Range("E35").Value = "  +5.10%  "
$ws.Range("E36").Value = "  +10.13%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").Value = "'6.47"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  +7.44%  "
$ws.Range("D40").Value = "'19.61"
$ws.Range("E40").Value = "  +9.95%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "'8.87"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.48"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'98.59"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +12.61%  "
$ws.Range("D49").Value = "'2.36"
$ws.Range("E49").Value = "  +5.57%  "
$ws.Range("D50").Value = "1.441.36"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +1.51%  "

# Clear the quote-prefix styling picked up from the forced-text entries
# above so cell styles stay identical to the original (style index 0).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
